# push start year to 2019
# The "IT" worksheet's Year/Initial value (cell B2) records the initial
# simulation start year. Update it from 2020 to 2019.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IT")
$ws.Range("B2").Value = 2019
